$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fStr = [string]$fVal
    if ($fStr.Length -ne 8) {
        continue
    }

    $year = [int]$fStr.Substring(0, 4)
    $month = [int]$fStr.Substring(4, 2)
    $day = [int]$fStr.Substring(6, 2)
    if ($month -lt 1 -or $month -gt 12 -or $day -lt 1 -or $day -gt 31) {
        continue
    }

    $newE = [int]$eVal - 1

    if ($newE -eq 0) {
        $eCell.Value = [int]$dVal
        $fCell.Value = 20251215
    } else {
        $eCell.Value = $newE
    }
}
